# Apply the changes described in the diff to the workbook.
$wb = $excel.ActiveWorkbook

# --- "About" sheet ---
$wsAbout = $wb.Worksheets.Item("About")
# C1 date value changed from 45320 to 45392 (serial date -> 2024-04-10)
$wsAbout.Range("C1").Value = 45392

# --- "MCF" sheet ---
$wsMCF = $wb.Worksheets.Item("MCF")

$wsMCF.Range("B2").Value = 1
$wsMCF.Range("B3").Value = 1
$wsMCF.Range("B4").Value = 1
$wsMCF.Range("B6").Value = 1
$wsMCF.Range("B10").Value = 1
$wsMCF.Range("B11").Value = 1
$wsMCF.Range("B12").Value = 1
$wsMCF.Range("B13").Value = 1
$wsMCF.Range("B14").Value = 1
$wsMCF.Range("B16").Value = 1
$wsMCF.Range("B17").Value = 1
$wsMCF.Range("B18").Value = 1

# Update selection on MCF sheet to B17, and make it the active sheet/tab
$wsMCF.Activate()
$wsMCF.Range("B17").Select()
